# [Modify] Change view for print of instruction set architecture before to return
#
# The commit swaps the "Feuil1" sheet from a normal-view, zoomed-in layout
# back to a Page Break Preview layout used for checking what prints, and
# widens the print area a few rows to cover rows up to 46, while also
# switching the page scaling from "fit to page" back to an explicit 74%
# print scale.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet          # "Feuil1" - already the selected/active sheet
$ws.Activate()

# --- Print area: Feuil1!$A$1:$M$42 -> Feuil1!$A$1:$M$46 --------------------
$ws.PageSetup.PrintArea = '$A$1:$M$46'

# --- Page scaling: go back to an explicit 74% scale (was "fit to page") ----
# Setting PageSetup.Zoom (rather than .Scale) switches the sheet out of the
# "fit to page" scaling mode and writes a plain percentage scale, matching
# <pageSetup .../ scale="74" .../> in the target workbook.
$ws.PageSetup.Zoom = 74

# --- Window / sheet view: switch to Page Break Preview, change zoom -------
$win = $excel.ActiveWindow
$win.View = 1        # xlNormalView   - capture the "normal view" zoom first
$win.Zoom = 115       # zoomScaleNormal after returning to normal view
$win.View = 2         # xlPageBreakPreview - the view used to check printing
$win.Zoom = 55        # zoomScale / zoomScaleSheetLayoutView while previewing

# --- Selection moves from O37 to K56 ---------------------------------------
$ws.Range("K56").Select() | Out-Null
